# Updates the cryptos list: refreshed prices / 1h volume %, and restores
# the correct SEI / BitcoinSV row order (rows 49-50).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free per-cell writes below. "Price" (column D) values are forced to
# text via a temporary "@" number format so strings like "1.00" or "309.04"
# are not reinterpreted as numbers by Excel's automatic type detection; the
# format/style is then restored so no stray style attributes are left behind.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.530.80'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.530.91'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.23%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.04'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.89'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.45%  '
$ws.Range('E7').Value = '  -1.20%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('E9').Value = '  -2.20%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.91'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0804'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.98%  '
$ws.Range('E12').Value = '  -1.47%  '
$ws.Range('E13').Value = '  +0.14%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.942.68'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.42%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.88'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.531.93'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.94%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.819'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.72%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.513.98'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.98%  '
$ws.Range('E19').Value = '  -0.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0953'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.18'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.65%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.18'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '242.86'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.89'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.17%  '
$ws.Range('E25').Value = '  -1.18%  '
$ws.Range('E27').Value = '  -3.42%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.30'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.37%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '39.19'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.86%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.12'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.59%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '155.96'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.76'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.78'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +14.10%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0795'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.31%  '
$ws.Range('E35').Value = '  -2.79%  '
$ws.Range('E36').Value = '  -4.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.27'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.99%  '
$ws.Range('E38').Value = '  -7.05%  '
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('E40').Value = '  +0.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.31'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +10.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.75'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.86%  '
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.30'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.47%  '
$ws.Range('E45').Value = '  -2.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.967.58'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.73%  '
$ws.Range('E47').Value = '  -1.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.770.44'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.19%  '
$ws.Range('B49').Value = 'BitcoinSV'
$ws.Range('C49').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '81.35'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.17%  '
$ws.Range('B50').Value = 'SEI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.866'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +12.28%  '
$ws.Range('E51').Value = '  -0.45%  '
